# Remove two account rows from the "Export" sheet:
#   - row with Conta 004399832 / Euler / 51086   (originally sheet row 5)
#   - row with Conta 005009026 / Edmur / 30      (originally sheet row 204)
#
# Delete from the bottom up so the earlier delete doesn't shift the row
# index of the later one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(204).Delete()
$ws.Rows(5).Delete()
